$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new confirmation-checkbox Q&A pair rows, following the same
# key/value/comment pattern used throughout the sheet.
$ws.Range("A32").Value = "pq_confirmation"
$ws.Range("B32").Value = "Please confirm that your answers are correct."

$ws.Range("A33").Value = "pq_confirm_answers"
$ws.Range("B33").Value = "Please read your answers carefully one more time and check this box if they are correct. Then you can proceed."

# Match the style/formatting used by the preceding row (row 31), which
# uses a different font than the default rows above it. Copy formats
# only so the existing style is reused instead of a new one created.
$ws.Range("A31:B31").Copy()
$ws.Range("A32:B33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B36").Select()
